$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.314.53'
$ws.Range("E2").Value = '  -3.46%  '

# Row 3
$ws.Range("D3").Value = '1.933.39'
$ws.Range("E3").Value = '  -3.73%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.00%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7192'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.47%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3292'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -8.08%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.92'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.98%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06909'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.03%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8034'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.38%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08064'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.61%  '

# Row 13
$ws.Range("D13").Value = '1.933.04'
$ws.Range("E13").Value = '  -3.75%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.411'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.61%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.75'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.16%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.52'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.70%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008408'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.17%  '

# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '30.307.15'
$ws.Range("E18").Value = '  -3.51%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '252.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -9.02%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.817'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.82%  '

# Row 21
$ws.Range("D21").Value = '2.182.98'
$ws.Range("E21").Value = '  -3.70%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.875'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.34%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.724'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.58%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.77%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.404'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.29%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.60%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1335'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -10.77%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.559'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.08%  '

# Row 31
$ws.Range("E31").Value = '  -1.41%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.404'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.83%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.188'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.75%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05118'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.19%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.225'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.17%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7407'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.06%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.748'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.30%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01969'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.16%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.830'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.21%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.608'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.23%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.87'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.22%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4461'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.02%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.996'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.23%  '

# Row 44
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8368'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.64%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.11'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.78%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.767'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.87%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.313'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.52%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.88%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05958'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.55%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.473'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.46%  '
